$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates: force text interpretation so numeric-looking
# strings like '1.00' / '4.01' are not coerced to numbers, then restore the
# default 'Normal' style so no stray per-cell formatting is introduced.
$priceCells = @("D2","D3","D5","D6","D7","D8","D9","D14","D17","D18","D19","D20","D21","D22","D24","D25","D32","D33","D34","D36","D37","D38","D42","D44","D45","D47","D50","D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value2 = '57.192.91'
$ws.Range("D3").Value2 = '3.071.00'
$ws.Range("D5").Value2 = '521.94'
$ws.Range("D6").Value2 = '135.62'
$ws.Range("D7").Value2 = '1.00'
$ws.Range("D8").Value2 = '3.067.93'
$ws.Range("D9").Value2 = '0.473'
$ws.Range("D14").Value2 = '3.596.94'
$ws.Range("D17").Value2 = '57.176.55'
$ws.Range("D18").Value2 = '3.054.05'
$ws.Range("D19").Value2 = '5.85'
$ws.Range("D20").Value2 = '12.40'
$ws.Range("D21").Value2 = '7.78'
$ws.Range("D22").Value2 = '347.50'
$ws.Range("D24").Value2 = '68.83'
$ws.Range("D25").Value2 = '0.498'
$ws.Range("D32").Value2 = '5.83'
$ws.Range("D33").Value2 = '20.97'
$ws.Range("D34").Value2 = '158.23'
$ws.Range("D36").Value2 = '1.11'
$ws.Range("D37").Value2 = '5.97'
$ws.Range("D38").Value2 = '25.35'
$ws.Range("D42").Value2 = '4.01'
$ws.Range("D44").Value2 = '2.410.78'
$ws.Range("D45").Value2 = '36.53'
$ws.Range("D47").Value2 = '3.106.93'
$ws.Range("D50").Value2 = '0.931'
$ws.Range("D51").Value2 = '19.28'
foreach ($addr in $priceCells) { $ws.Range($addr).Style = "Normal" }

# Volume(1h) (column E) updates: plain text, Excel leaves these as strings
# on their own because of the leading/trailing spaces and '%' sign.
$ws.Range("E2").Value2 = '  -2.14%  '
$ws.Range("E3").Value2 = '  -2.15%  '
$ws.Range("E4").Value2 = '  -0.06%  '
$ws.Range("E5").Value2 = '  -2.25%  '
$ws.Range("E6").Value2 = '  -5.26%  '
$ws.Range("E7").Value2 = '  +0.00%  '
$ws.Range("E8").Value2 = '  -2.31%  '
$ws.Range("E9").Value2 = '  +4.95%  '
$ws.Range("E10").Value2 = '  +0.71%  '
$ws.Range("E11").Value2 = '  -3.26%  '
$ws.Range("E12").Value2 = '  +1.23%  '
$ws.Range("E13").Value2 = '  +1.48%  '
$ws.Range("E14").Value2 = '  -2.27%  '
$ws.Range("E15").Value2 = '  -2.28%  '
$ws.Range("E16").Value2 = '  -4.49%  '
$ws.Range("E17").Value2 = '  -2.30%  '
$ws.Range("E18").Value2 = '  -2.81%  '
$ws.Range("E19").Value2 = '  -4.40%  '
$ws.Range("E20").Value2 = '  -3.85%  '
$ws.Range("E21").Value2 = '  -2.65%  '
$ws.Range("E22").Value2 = '  +1.44%  '
$ws.Range("E23").Value2 = '  -0.17%  '
$ws.Range("E24").Value2 = '  +1.52%  '
$ws.Range("E25").Value2 = '  -3.08%  '
$ws.Range("E26").Value2 = '  +0.04%  '
$ws.Range("E27").Value2 = '  -3.11%  '
$ws.Range("E28").Value2 = '  -9.82%  '
$ws.Range("E29").Value2 = '  -0.01%  '
$ws.Range("E30").Value2 = '  -5.53%  '
$ws.Range("E31").Value2 = '  -3.34%  '
$ws.Range("E32").Value2 = '  -10.04%  '
$ws.Range("E33").Value2 = '  -0.91%  '
$ws.Range("E34").Value2 = '  -0.12%  '
$ws.Range("E35").Value2 = '  -0.19%  '
$ws.Range("E36").Value2 = '  -7.14%  '
$ws.Range("E37").Value2 = '  -4.59%  '
$ws.Range("E38").Value2 = '  -3.13%  '
$ws.Range("E39").Value2 = '  -3.82%  '
$ws.Range("E40").Value2 = '  -2.77%  '
$ws.Range("E41").Value2 = '  -6.18%  '
$ws.Range("E42").Value2 = '  -0.31%  '
$ws.Range("E43").Value2 = '  -2.87%  '
$ws.Range("E44").Value2 = '  +4.66%  '
$ws.Range("E45").Value2 = '  -0.27%  '
$ws.Range("E46").Value2 = '  -0.10%  '
$ws.Range("E47").Value2 = '  -2.32%  '
$ws.Range("E48").Value2 = '  -2.49%  '
$ws.Range("E49").Value2 = '  -2.43%  '
$ws.Range("E50").Value2 = '  -8.14%  '
$ws.Range("E51").Value2 = '  -6.80%  '
